$wb = $excel.ActiveWorkbook

# "Ready for handoff" no longer applies to any file in this report snapshot;
# clear the per-locale status cells on the Overview sheet and the Status
# column on each locale sheet (leaving them as empty text, matching the
# "In Translation" cells already in those columns).
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E5:F7").Value = "'"
$overview.Range("E5:F7").Style = "Normal"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C5:C7").Value = "'"
$zhcn.Range("C5:C7").Style = "Normal"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C5:C7").Value = "'"
$dede.Range("C5:C7").Style = "Normal"

# Now that the column no longer contains the long "Ready for handoff"
# string, autofit the affected columns so they shrink to fit what's left.
$overview.Columns.Item(5).EntireColumn.AutoFit()
$overview.Columns.Item(6).EntireColumn.AutoFit()
$zhcn.Columns.Item(3).EntireColumn.AutoFit()
$dede.Columns.Item(3).EntireColumn.AutoFit()
